$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 728.7273
$ws.Range("I33").Value = 688.4286
$ws.Range("K33").Value = 688.4286
$ws.Range("M33").Value = -459.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 20919.285
$ws.Range("I70").Value = 1858
$ws.Range("J70").Value = 46334.332
$ws.Range("K70").Value = 5574
$ws.Range("L70").Value = 139002.996
$ws.Range("M70").Value = -5304
$ws.Range("N70").Value = -139542.996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 20919.285
$ws.Range("I73").Value = 1858
$ws.Range("J73").Value = 46334.332
$ws.Range("K73").Value = 5574
$ws.Range("L73").Value = 139002.996
$ws.Range("M73").Value = -4638
$ws.Range("N73").Value = -140874.996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 633
$ws.Range("I80").Value = 564.9231
$ws.Range("K80").Value = 1694.7693
$ws.Range("M80").Value = -696.7692999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 633
$ws.Range("I83").Value = 564.9231
$ws.Range("K83").Value = 5084.3079
$ws.Range("M83").Value = -92.30789999999979

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 11192.667
$ws.Range("J116").Value = 13056
$ws.Range("L116").Value = 13056
$ws.Range("N116").Value = -19940

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1516359.5
$ws.Range("I137").Value = 869220.4
$ws.Range("K137").Value = 2607661.2
$ws.Range("M137").Value = -2605111.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3199.9033
$ws.Range("I138").Value = 2033.7142
$ws.Range("J138").Value = 5648.9
$ws.Range("K138").Value = 6101.142599999999
$ws.Range("L138").Value = 16946.7
$ws.Range("M138").Value = -961.1425999999992
$ws.Range("N138").Value = -27226.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10423164
$ws.Range("I32").Value = 11115331
$ws.Range("K32").Value = 11115331
$ws.Range("M32").Value = -11115044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2794170.2
$ws.Range("I61").Value = 3352162.5
$ws.Range("K61").Value = 3352162.5
$ws.Range("M61").Value = -3351950.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3381646.8
$ws.Range("I74").Value = 4034810
$ws.Range("K74").Value = 4034810
$ws.Range("M74").Value = -4033936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3381646.8
$ws.Range("I77").Value = 4034810
$ws.Range("K77").Value = 20174050
$ws.Range("M77").Value = -20169682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 12449.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 12449.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 37348.5
$ws.Range("N132").Value = -42408.5
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 87466.86
$ws.Range("J134").Value = 87466.86
$ws.Range("L134").Value = 87466.86
$ws.Range("N134").Value = -97606.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2794170.2
$ws.Range("I136").Value = 3352162.5
$ws.Range("K136").Value = 10056487.5
$ws.Range("M136").Value = -10053937.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1599.1538
$ws.Range("I20").Value = 1797.7646
$ws.Range("J20").Value = 1224
$ws.Range("K20").Value = 1797.7646
$ws.Range("L20").Value = 1224
$ws.Range("M20").Value = -1550.7646
$ws.Range("N20").Value = -1718

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 32974
$ws.Range("J88").Value = 32974
$ws.Range("L88").Value = 32974
$ws.Range("N88").Value = -33786

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 32974
$ws.Range("J91").Value = 32974
$ws.Range("L91").Value = 32974
$ws.Range("N91").Value = -35782

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 16500
$ws.Range("I96").Value = 16500
$ws.Range("K96").Value = 16500
$ws.Range("M96").Value = -13754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1037789.06
$ws.Range("J134").Value = 13973.75
$ws.Range("L134").Value = 41921.25
$ws.Range("N134").Value = -46991.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 99984
$ws.Range("J135").Value = 99984
$ws.Range("L135").Value = 99984
$ws.Range("N135").Value = -110124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 212.46153
$ws.Range("J7").Value = 419.5
$ws.Range("L7").Value = 419.5
$ws.Range("N7").Value = -645.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8347156.5
$ws.Range("I132").Value = 15688.23
$ws.Range("K132").Value = 47064.69
$ws.Range("M132").Value = -44534.69

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1813.3158
$ws.Range("I134").Value = 1864.3889
$ws.Range("K134").Value = 5593.1667
$ws.Range("M134").Value = -3058.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 90793.71000000001
$ws.Range("J140").Value = 90793.71000000001
$ws.Range("L140").Value = 90793.71000000001
$ws.Range("N140").Value = -101153.71

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 134500180
$ws.Range("I4").Value = 163538670
$ws.Range("K4").Value = 490616010
$ws.Range("M4").Value = -490615898

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 153.16667
$ws.Range("I17").Value = 70.333336
$ws.Range("J17").Value = 733
$ws.Range("K17").Value = 211.000008
$ws.Range("L17").Value = 2199
$ws.Range("M17").Value = -42.00000800000001
$ws.Range("N17").Value = -2537

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 678.8
$ws.Range("I23").Value = 115
$ws.Range("J23").Value = 819.75
$ws.Range("K23").Value = 345
$ws.Range("L23").Value = 2459.25
$ws.Range("M23").Value = -110
$ws.Range("N23").Value = -2929.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2124.7693
$ws.Range("J113").Value = 2097
$ws.Range("L113").Value = 6291
$ws.Range("N113").Value = -10631

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10184.519
$ws.Range("J131").Value = 10502
$ws.Range("L131").Value = 31506
$ws.Range("N131").Value = -41586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1723236.9
$ws.Range("I132").Value = 2009610.6
$ws.Range("K132").Value = 6028831.800000001
$ws.Range("M132").Value = -6026301.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6697.9414
$ws.Range("I136").Value = 5633.2856
$ws.Range("K136").Value = 16899.8568
$ws.Range("M136").Value = -14349.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3374
$ws.Range("I113").Value = 1805.3572
$ws.Range("J113").Value = 5204.0835
$ws.Range("K113").Value = 5416.071599999999
$ws.Range("L113").Value = 15612.2505
$ws.Range("M113").Value = -3246.071599999999
$ws.Range("N113").Value = -19952.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9530593
$ws.Range("I136").Value = 10891263
$ws.Range("J136").Value = 5899.8
$ws.Range("K136").Value = 32673789
$ws.Range("L136").Value = 17699.4
$ws.Range("N136").Value = -22799.4
$ws.Range("M136").Value = -32671239
